$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column C ("cantidad") entirely - columns shift left
$ws.Columns.Item(3).Delete()

# The "precio_venta" column (now column C) no longer needs the Text ("@") number format
$ws.Range("C2:C6").NumberFormat = "General"

# Center-align (horizontally) all data rows in the remaining columns A:E
$ws.Range("A2:E6").HorizontalAlignment = -4108
